# Slide 1 (cId 2474004964 / sldId 256): the course-code title placeholder
# ("Title 1", ctrTitle) is updated from "MHW1" to "MHW3".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Title
$titleShape.TextFrame.TextRange.Text = "MHW3"
